# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @{ B = 0.3048080303191223; C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 2.626907116734944 }
    3 = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
    4 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489; E = 0.496779210170732;  G = 9.295990156953671 }
    5 = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Range("B$row").Value = $rowData.B
    $ws.Range("C$row").Value = $rowData.C
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("E$row").Value = $rowData.E
    $ws.Range("G$row").Value = $rowData.G
}
